$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking snapshot refresh: Price (D) and Volume(1h) (E) columns are
# stored as text in this sheet (trailing zeros / thousands-dot formatting
# must survive, e.g. "1.000" or "30.324.91"), so force text for any new
# value that Excel would otherwise auto-convert to a number, then drop the
# quote-prefix style Excel applies so the cell format matches the original.
function Set-TextValue {
    param($range, [string]$text)
    if ($text -match '^[+-]?\d+(\.\d+)?([eE][+-]?\d+)?$') {
        $range.Value = "'" + $text
        $range.Style = "Normal"
    } else {
        $range.Value = $text
    }
}

Set-TextValue $ws.Range("D2") '30.324.91'
Set-TextValue $ws.Range("E2") '  -3.10%  '
Set-TextValue $ws.Range("D3") '1.936.61'
Set-TextValue $ws.Range("E3") '  -3.15%  '
Set-TextValue $ws.Range("D4") '1.000'
Set-TextValue $ws.Range("E4") '  +0.25%  '
Set-TextValue $ws.Range("D5") '250.65'
Set-TextValue $ws.Range("E5") '  -1.80%  '
Set-TextValue $ws.Range("D6") '0.7234'
Set-TextValue $ws.Range("E6") '  -7.40%  '
Set-TextValue $ws.Range("D7") '1.000'
Set-TextValue $ws.Range("E7") '  +0.20%  '
Set-TextValue $ws.Range("D8") '0.3316'
Set-TextValue $ws.Range("E8") '  -4.90%  '
Set-TextValue $ws.Range("D9") '27.93'
Set-TextValue $ws.Range("E9") '  -0.96%  '
Set-TextValue $ws.Range("D10") '0.07259'
Set-TextValue $ws.Range("E10") '  +2.27%  '
Set-TextValue $ws.Range("D11") '0.8109'
Set-TextValue $ws.Range("E11") '  -3.97%  '
Set-TextValue $ws.Range("D12") '0.08086'
Set-TextValue $ws.Range("E12") '  -1.40%  '
Set-TextValue $ws.Range("D13") '1.939.74'
Set-TextValue $ws.Range("E13") '  -2.92%  '
Set-TextValue $ws.Range("D14") '5.509'
Set-TextValue $ws.Range("E14") '  -2.62%  '
Set-TextValue $ws.Range("D15") '94.79'
Set-TextValue $ws.Range("E15") '  -6.24%  '
Set-TextValue $ws.Range("D16") '15.13'
Set-TextValue $ws.Range("E16") '  -1.65%  '
Set-TextValue $ws.Range("D17") '30.340.46'
Set-TextValue $ws.Range("E17") '  -3.02%  '
Set-TextValue $ws.Range("D18") '0.000008304'
Set-TextValue $ws.Range("E18") '  +3.03%  '
Set-TextValue $ws.Range("D19") '254.11'
Set-TextValue $ws.Range("E19") '  -7.11%  '
Set-TextValue $ws.Range("D20") '5.900'
Set-TextValue $ws.Range("E20") '  -1.61%  '
Set-TextValue $ws.Range("D21") '2.193.54'
Set-TextValue $ws.Range("E21") '  -2.92%  '
Set-TextValue $ws.Range("D22") '1.000'
Set-TextValue $ws.Range("E22") '  +0.18%  '
Set-TextValue $ws.Range("D23") '1.000'
Set-TextValue $ws.Range("E23") '  +0.26%  '
Set-TextValue $ws.Range("D24") '6.986'
Set-TextValue $ws.Range("E24") '  -1.93%  '
Set-TextValue $ws.Range("D25") '9.772'
Set-TextValue $ws.Range("E25") '  -3.20%  '
Set-TextValue $ws.Range("D26") '163.58'
Set-TextValue $ws.Range("E26") '  -0.63%  '
Set-TextValue $ws.Range("D27") '2.395'
Set-TextValue $ws.Range("E27") '  -0.95%  '
Set-TextValue $ws.Range("D28") '19.31'
Set-TextValue $ws.Range("E28") '  -3.43%  '
Set-TextValue $ws.Range("D29") '0.1318'
Set-TextValue $ws.Range("E29") '  -7.93%  '
Set-TextValue $ws.Range("D30") '1.570'
Set-TextValue $ws.Range("E30") '  -1.71%  '
Set-TextValue $ws.Range("D31") '1.351'
Set-TextValue $ws.Range("E31") '  -1.05%  '
Set-TextValue $ws.Range("D32") '4.445'
Set-TextValue $ws.Range("E32") '  -4.48%  '
Set-TextValue $ws.Range("D33") '4.179'
Set-TextValue $ws.Range("E33") '  -6.23%  '
Set-TextValue $ws.Range("D34") '0.05214'
Set-TextValue $ws.Range("E34") '  -2.82%  '
Set-TextValue $ws.Range("D35") '1.290'
Set-TextValue $ws.Range("E35") '  +1.56%  '
Set-TextValue $ws.Range("D36") '0.7514'
Set-TextValue $ws.Range("E36") '  -5.06%  '
Set-TextValue $ws.Range("E37") '  -0.86%  '
Set-TextValue $ws.Range("D38") '0.01986'
Set-TextValue $ws.Range("E38") '  -1.14%  '
Set-TextValue $ws.Range("D39") '2.823'
Set-TextValue $ws.Range("E39") '  -3.24%  '
Set-TextValue $ws.Range("D40") '79.57'
Set-TextValue $ws.Range("E40") '  -6.83%  '
Set-TextValue $ws.Range("E41") '  -6.96%  '
Set-TextValue $ws.Range("D42") '0.4553'
Set-TextValue $ws.Range("E42") '  -3.11%  '
Set-TextValue $ws.Range("D43") '2.032'
Set-TextValue $ws.Range("E43") '  -5.91%  '
Set-TextValue $ws.Range("D44") '0.8470'
Set-TextValue $ws.Range("E44") '  -1.32%  '
Set-TextValue $ws.Range("D45") '0.9999'
Set-TextValue $ws.Range("E45") '  +0.15%  '
Set-TextValue $ws.Range("D46") '101.91'
Set-TextValue $ws.Range("E46") '  -3.54%  '
Set-TextValue $ws.Range("D47") '9.741'
Set-TextValue $ws.Range("E47") '  -5.72%  '
Set-TextValue $ws.Range("D48") '7.484'
Set-TextValue $ws.Range("E48") '  -3.86%  '
Set-TextValue $ws.Range("D49") '36.87'
Set-TextValue $ws.Range("E49") '  -2.67%  '
Set-TextValue $ws.Range("D50") '0.4195'
Set-TextValue $ws.Range("E50") '  -3.52%  '
Set-TextValue $ws.Range("D51") '0.06036'
Set-TextValue $ws.Range("E51") '  +0.05%  '
